$p = $ppt.ActivePresentation
$s = $p.Slides.Item(6)

# Locate the content placeholder that lists the simulation parameters
# (Data Rate / Delay / Packet Size / Max Packets / Simulation Time).
$shape = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $candidate = $s.Shapes.Item($i)
    if ($candidate.HasTextFrame -and $candidate.TextFrame.TextRange.Text -like "*Packet Size*") {
        $shape = $candidate
    }
}

$tf = $shape.TextFrame
$tr = $tf.TextRange

# "  Packet Size = 1024" -> "  Packet Size = 1024 bytes"
$packetSizePara = $tr.Paragraphs(6)
$packetSizePara.Runs(1).Text = "  Packet Size = 1024 bytes"

# "  Simulation Time = 11.0s" -> split into two runs:
#   "  Simulation Time " / "= 11.0s"
$simTimePara = $tr.Paragraphs(8)
$simTimeRun = $simTimePara.Runs(1)
$simTimeRun.Text = "  Simulation Time "
$simTimeRun.InsertAfter("= 11.0s") | Out-Null
